$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp values for rows 2..93 (column F), in order
$timestamps = @(
    "2021-10-05 10:50:26.603353",
    "2021-10-05 10:50:26.603367",
    "2021-10-05 10:50:26.603371",
    "2021-10-05 10:50:26.603374",
    "2021-10-05 10:50:26.603380",
    "2021-10-05 10:50:26.603386",
    "2021-10-05 10:50:26.603390",
    "2021-10-05 10:50:26.603393",
    "2021-10-05 10:50:26.603397",
    "2021-10-05 10:50:26.603400",
    "2021-10-05 10:50:26.603403",
    "2021-10-05 10:50:26.603406",
    "2021-10-05 10:50:26.603409",
    "2021-10-05 10:50:26.603412",
    "2021-10-05 10:50:26.603416",
    "2021-10-05 10:50:26.603419",
    "2021-10-05 10:50:26.603422",
    "2021-10-05 10:50:26.603426",
    "2021-10-05 10:50:26.603429",
    "2021-10-05 10:50:26.603432",
    "2021-10-05 10:50:26.603435",
    "2021-10-05 10:50:26.603439",
    "2021-10-05 10:50:26.603442",
    "2021-10-05 10:50:26.603445",
    "2021-10-05 10:50:26.603449",
    "2021-10-05 10:50:26.603453",
    "2021-10-05 10:50:26.603456",
    "2021-10-05 10:50:26.603459",
    "2021-10-05 10:50:26.603462",
    "2021-10-05 10:50:26.603465",
    "2021-10-05 10:50:26.603468",
    "2021-10-05 10:50:26.603471",
    "2021-10-05 10:50:26.603475",
    "2021-10-05 10:50:26.603479",
    "2021-10-05 10:50:26.603482",
    "2021-10-05 10:50:26.603485",
    "2021-10-05 10:50:26.603488",
    "2021-10-05 10:50:26.603491",
    "2021-10-05 10:50:26.603494",
    "2021-10-05 10:50:26.603497",
    "2021-10-05 10:50:26.603500",
    "2021-10-05 10:50:26.603503",
    "2021-10-05 10:50:26.603506",
    "2021-10-05 10:50:26.603509",
    "2021-10-05 10:50:26.603512",
    "2021-10-05 10:50:26.603516",
    "2021-10-05 10:50:26.603518",
    "2021-10-05 10:50:26.603522",
    "2021-10-05 10:50:26.603525",
    "2021-10-05 10:50:26.603528",
    "2021-10-05 10:50:26.603531",
    "2021-10-05 10:50:26.603534",
    "2021-10-05 10:50:26.603537",
    "2021-10-05 10:50:26.603540",
    "2021-10-05 10:50:26.603543",
    "2021-10-05 10:50:26.603546",
    "2021-10-05 10:50:26.603549",
    "2021-10-05 10:50:26.603553",
    "2021-10-05 10:50:26.603556",
    "2021-10-05 10:50:26.603559",
    "2021-10-05 10:50:26.603562",
    "2021-10-05 10:50:26.603565",
    "2021-10-05 10:50:26.603568",
    "2021-10-05 10:50:26.603571",
    "2021-10-05 10:50:26.603576",
    "2021-10-05 10:50:26.603579",
    "2021-10-05 10:50:26.603582",
    "2021-10-05 10:50:26.603585",
    "2021-10-05 10:50:26.603588",
    "2021-10-05 10:50:26.603592",
    "2021-10-05 10:50:26.603595",
    "2021-10-05 10:50:26.603598",
    "2021-10-05 10:50:26.603601",
    "2021-10-05 10:50:26.603604",
    "2021-10-05 10:50:26.603607",
    "2021-10-05 10:50:26.603611",
    "2021-10-05 10:50:26.603616",
    "2021-10-05 10:50:26.603620",
    "2021-10-05 10:50:26.603623",
    "2021-10-05 10:50:26.603626",
    "2021-10-05 10:50:26.603629",
    "2021-10-05 10:50:26.603632",
    "2021-10-05 10:50:26.603635",
    "2021-10-05 10:50:26.603639",
    "2021-10-05 10:50:26.603642",
    "2021-10-05 10:50:26.603645",
    "2021-10-05 10:50:26.603648",
    "2021-10-05 10:50:26.603651",
    "2021-10-05 10:50:26.603654",
    "2021-10-05 10:50:26.603657",
    "2021-10-05 10:50:26.603660",
    "2021-10-05 10:50:26.603663"
)

# Header cell F1: copy formatting from E1 (bold/centered header style) then set its text
$ws.Cells.Item(1, 5).Copy()
$ws.Cells.Item(1, 6).PasteSpecial(-4122)
$ws.Cells.Item(1, 6).Value = "time_taken"

# Fill F2:F93 with the recorded time_taken values (plain, unstyled like the other data cells)
for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
